$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Solicitud gráfica")
try {
    $ws.Unprotect()
    Write-Output "unprotect ok no pw"
} catch {
    Write-Output "ERROR unprotect: $_"
}
